# Update the Handoff/Handback datetime stamps for the 85f86cf4... file row
# on both the zh-cn and de-de report sheets, to reflect a newer report
# generation run (commit: "Generate Report for Handback").

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-19 02:35:01"
$wsZh.Range("H4").Value = "2016-03-19 02:35:26"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-19 02:35:05"
$wsDe.Range("H4").Value = "2016-03-19 02:35:31"
